# Fixed issue #21: the "medida" metadata block in row 3/4/5, column I
# (Direccion Provincial Nombre) was incorrectly generated as a measure
# (iaest-measure:direccion-provincial-nombre / medida / xsd:string)
# instead of a dimension that references the area, like columns D/F
# (sdmx-dimension:refArea / dim / URI-Provincia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "medida"/"dim" identifier for column I -> sdmx-dimension:refArea
$i3 = $ws.Cells.Item(3, 9)
$i3.Value = "sdmx-dimension:refArea"

# Give I3 a distinct look (black Arial text on a solid white fill) so the
# corrected dimension column stands out, same as the source workbook.
$i3.Font.Name = "Arial"
$i3.Font.Color = 0
$i3.Interior.PatternColor = 16777215
$i3.Interior.Color = 16777215
$i3.Interior.Pattern = 1

# Row 4: whether the column is a "medida" (measure) or a "dim" (dimension)
$ws.Cells.Item(4, 9).Value = "dim"

# Row 5: the datatype/URI associated with column I
$ws.Cells.Item(5, 9).Value = "URI-Provincia"
